# Slide 18 "Example: Subprogram Inlining" - reduce the "space before" on
# the four code-block paragraphs (the "proc inc(...)", "{", "x := x + 1;",
# "}" lines) from 3 points (300 hundredths) to 1 point (100 hundredths).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tr = $sh.TextFrame.TextRange

for ($i = 1; $i -le 4; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.ParagraphFormat.SpaceBefore = 1
}
